$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.195.54"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "2.413.59"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.78"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.20"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.148"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.350"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.79"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "2.840.36"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "60.069.49"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "2.411.31"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.20"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.52"
$ws.Range("E19").Value = "  +3.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "328.32"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.79"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.81"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.175"
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.60"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.81"
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("D29").Value = "0.0₃0771"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.20"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  +8.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.402"
$ws.Range("E33").Value = "  -2.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.46"
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.21"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "324.78"
$ws.Range("E39").Value = "  +3.41%  "
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.44"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "146.40"
$ws.Range("E42").Value = "  +5.66%  "
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0972"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.88"
$ws.Range("E45").Value = "  +1.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0517"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.577"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0222"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  -0.81%  "
